$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data between row 2 and row 4 for columns D, J, K, L, M, O, P

# Save row 2 current values
$d2 = $ws.Range("D2").Value()
$j2 = $ws.Range("J2").Value()
$k2 = $ws.Range("K2").Value()
$l2 = $ws.Range("L2").Value()
$m2 = $ws.Range("M2").Value()
$o2 = $ws.Range("O2").Value()
$p2 = $ws.Range("P2").Value()

# Save row 4 current values
$d4 = $ws.Range("D4").Value()
$j4 = $ws.Range("J4").Value()
$k4 = $ws.Range("K4").Value()
$l4 = $ws.Range("L4").Value()
$m4 = $ws.Range("M4").Value()
$o4 = $ws.Range("O4").Value()
$p4 = $ws.Range("P4").Value()

# Write row 4 values into row 2
$ws.Range("D2").Value = $d4
$ws.Range("J2").Value = $j4
$ws.Range("K2").Value = $k4
$ws.Range("L2").Value = $l4
$ws.Range("M2").Value = $m4
$ws.Range("O2").Value = $o4
$ws.Range("P2").Value = $p4

# Write row 2 (original) values into row 4
$ws.Range("D4").Value = $d2
$ws.Range("J4").Value = $j2
$ws.Range("K4").Value = $k2
$ws.Range("L4").Value = $l2
$ws.Range("M4").Value = $m2
$ws.Range("O4").Value = $o2
$ws.Range("P4").Value = $p2
